$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default, unstyled) taken from a data cell (D4) that keeps its General/text
# formatting unmodified, used to avoid introducing spurious style indices when we
# temporarily switch a cell to Text format so Excel does not coerce numeric-looking
# strings (e.g. "58.50", "1.00") into actual numbers.
$defaultStyle = $ws.Range("D4").Style

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.780.68"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "  -0.61%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.219.26"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "  -0.41%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.36"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +4.35%  "

# Row 6
$ws.Range("E6").Value = "  +0.68%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.35"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  +3.71%  "

# Row 8
$ws.Range("E8").Value = "  +0.13%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.599"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  +9.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.38"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +11.08%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0967"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  +0.18%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.50"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +0.78%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.28"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  +9.05%  "

# Row 14
$ws.Range("E14").Value = "  +0.11%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.551.26"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  -0.34%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.903"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +6.00%  "

# Row 17
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "15.02"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +1.56%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.222.51"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  +0.14%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.747.65"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  -0.19%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0966"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +1.22%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.29"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  +2.81%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.58"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.14"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  +0.41%  "

# Row 24
$ws.Range("E24").Value = "  +0.42%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.05"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  +12.75%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.05"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +22.30%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +0.06%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.54"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  +4.43%  "

# Row 29
$ws.Range("E29").Value = "  -1.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.08"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +0.40%  "

# Row 31
$ws.Range("E31").Value = "  +2.34%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.124"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +2.55%  "

# Row 33
$ws.Range("E33").Value = "  +9.24%  "

# Row 34
$ws.Range("E34").Value = "  -0.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0747"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  +4.97%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.71"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +0.95%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.03"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  +12.20%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.02"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +4.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0304"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +8.75%  "

# Row 40
$ws.Range("E40").Value = "  -0.49%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.93"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +2.24%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.32"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +26.50%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.05"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -2.45%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.204"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  +7.75%  "

# Row 45
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.95"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +1.08%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.86"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -1.80%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.102"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  +2.00%  "

# Row 48
$ws.Range("E48").Value = "  +0.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.54"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  -0.37%  "

# Row 50
$ws.Range("E50").Value = "  +6.48%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.42"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  +5.66%  "
